$wb = $excel.ActiveWorkbook

# The new "Turkey" sheet reuses the exact same layout/styles as the
# existing "Spain" sheet, so build it by duplicating "Spain" and placing
# the copy right after it (i.e. at the end of the workbook).
$spain = $wb.Worksheets.Item("Spain")
$spain.Copy($null, $spain)

$turkey = $wb.Worksheets.Item($spain.Index + 1)
$turkey.Name = "Turkey"

# Fill in the market-specific cells.
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3325/T3324"

# Re-fit row heights (3-5 no longer need the taller wrapped height) and
# column B's width for the new content.
$turkey.Rows.Item(3).AutoFit()
$turkey.Rows.Item(4).AutoFit()
$turkey.Rows.Item(5).AutoFit()
$turkey.Columns.Item(2).AutoFit()

# "Spain" is no longer the active tab; its selection collapses to the
# whole used range.
$spain.Activate()
$spain.Range("A1:D10").Select() | Out-Null

# "Turkey" becomes the new active sheet/tab.
$turkey.Activate()
$turkey.Range("E19").Select() | Out-Null
